# Update the "想去人数" (interested-count) column F for the two sheets that
# carry the full event list: "展览" and "全部类型". The other two sheets
# ("演出", "本地生活") only contain a header row and are left untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 150
    3  = 1696
    4  = 785
    5  = 1118
    7  = 11877
    9  = 97
    11 = 404
    12 = 1107
    14 = 13454
    15 = 13388
    19 = 37
    20 = 275
    23 = 158
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
